$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert cells (shift right) only within the header/data/total rows (8-10), column F.
$ws.Range("F8:F10").Insert(-4161)

# Header for the newly inserted column
$ws.Cells.Item(8, 6).Value = "Форма зайнятості"

$ws.Range("Q17").Select()
